$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 213.16667
$ws.Range("I6").Value = 242
$ws.Range("K6").Value = 726
$ws.Range("M6").Value = -614
$ws.Range("H9").Value = 668.9091
$ws.Range("I9").Value = 754.55554
$ws.Range("J9").Value = 283.5
$ws.Range("K9").Value = 754.55554
$ws.Range("L9").Value = 283.5
$ws.Range("M9").Value = -585.55554
$ws.Range("N9").Value = -621.5
$ws.Range("H43").Value = 5674.1665
$ws.Range("I43").Value = 4995.5
$ws.Range("K43").Value = 4995.5
$ws.Range("M43").Value = -4926.5
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H64").Value = 7482.8335
$ws.Range("I64").Value = 6631.6665
$ws.Range("J64").Value = 8334
$ws.Range("K64").Value = 6631.6665
$ws.Range("L64").Value = 8334
$ws.Range("M64").Value = -6383.6665
$ws.Range("N64").Value = -8830
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H67").Value = 7482.8335
$ws.Range("I67").Value = 6631.6665
$ws.Range("J67").Value = 8334
$ws.Range("K67").Value = 6631.6665
$ws.Range("L67").Value = 8334
$ws.Range("M67").Value = -5773.6665
$ws.Range("N67").Value = -10050
$ws.Range("H87").Value = 86724.836
$ws.Range("J87").Value = 86724.836
$ws.Range("L87").Value = 86724.836
$ws.Range("N87").Value = -89220.836
$ws.Range("H90").Value = 86724.836
$ws.Range("J90").Value = 86724.836
$ws.Range("L90").Value = 260174.508
$ws.Range("N90").Value = -272654.508
$ws.Range("H111").Value = 3000
$ws.Range("I111").Value = 3000
$ws.Range("K111").Value = 9000
$ws.Range("M111").Value = -5933
$ws.Range("H117").Value = 34999.5
$ws.Range("J117").Value = 34999.5
$ws.Range("L117").Value = 34999.5
$ws.Range("N117").Value = -44177.5
$ws.Range("H138").Value = 1976.4884
$ws.Range("J138").Value = 2129.7585
$ws.Range("L138").Value = 6389.2755
$ws.Range("N138").Value = -16669.2755

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H92").Value = 27407.334
$ws.Range("J92").Value = 27407.334
$ws.Range("L92").Value = 27407.334
$ws.Range("N92").Value = -32399.334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 11917.333
$ws.Range("I26").Value = 11917.333
$ws.Range("K26").Value = 11917.333
$ws.Range("M26").Value = -11625.333
$ws.Range("H86").Value = 8341
$ws.Range("I86").Value = 7567
$ws.Range("K86").Value = 7567
$ws.Range("M86").Value = -6444
$ws.Range("H89").Value = 8341
$ws.Range("I89").Value = 7567
$ws.Range("K89").Value = 37835
$ws.Range("M89").Value = -32219
$ws.Range("H138").Value = 147500
$ws.Range("I138").Value = 135000
$ws.Range("K138").Value = 135000
$ws.Range("M138").Value = -129860

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 250
$ws.Range("I14").Value = 250
$ws.Range("K14").Value = 250
$ws.Range("M14").Value = -80
$ws.Range("H20").Value = 30000
$ws.Range("J20").Value = 30000
$ws.Range("L20").Value = 30000
$ws.Range("N20").Value = -30472
$ws.Range("H28").Value = 11928
$ws.Range("J28").Value = 11928
$ws.Range("L28").Value = 11928
$ws.Range("N28").Value = -12418
$ws.Range("H30").Value = 30000
$ws.Range("J30").Value = 30000
$ws.Range("L30").Value = 30000
$ws.Range("N30").Value = -30182
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()
$ws.Range("H63").Value = 100271
$ws.Range("J63").Value = 100271
$ws.Range("L63").Value = 100271
$ws.Range("N63").Value = -101643
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H66").Value = 100271
$ws.Range("J66").Value = 100271
$ws.Range("L66").Value = 300813
$ws.Range("N66").Value = -307677
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H69").Value = 9666.666999999999
$ws.Range("I69").Value = 10000
$ws.Range("J69").Value = 9000
$ws.Range("K69").Value = 10000
$ws.Range("L69").Value = 9000
$ws.Range("M69").Value = -9251
$ws.Range("N69").Value = -10498
$ws.Range("H72").Value = 9666.666999999999
$ws.Range("I72").Value = 10000
$ws.Range("J72").Value = 9000
$ws.Range("K72").Value = 30000
$ws.Range("L72").Value = 27000
$ws.Range("M72").Value = -26256
$ws.Range("N72").Value = -34488
$ws.Range("H75").Value = 39973.332
$ws.Range("J75").Value = 39973.332
$ws.Range("L75").Value = 39973.332
$ws.Range("N75").Value = -41969.332
$ws.Range("H78").Value = 39973.332
$ws.Range("J78").Value = 39973.332
$ws.Range("L78").Value = 119919.996
$ws.Range("N78").Value = -129903.996
$ws.Range("H81").Value = 39999.5
$ws.Range("J81").Value = 39999.5
$ws.Range("L81").Value = 39999.5
$ws.Range("N81").Value = -41995.5
$ws.Range("H84").Value = 39999.5
$ws.Range("J84").Value = 39999.5
$ws.Range("L84").Value = 119998.5
$ws.Range("N84").Value = -129982.5
$ws.Range("H86").Value = 4258.75
$ws.Range("I86").Value = 4100.4546
$ws.Range("J86").Value = 6000
$ws.Range("K86").Value = 4100.4546
$ws.Range("L86").Value = 6000
$ws.Range("M86").Value = -2977.4546
$ws.Range("N86").Value = -8246
$ws.Range("H88").Value = 6000
$ws.Range("J88").Value = 5000
$ws.Range("L88").Value = 5000
$ws.Range("N88").Value = -5812
$ws.Range("H89").Value = 4258.75
$ws.Range("I89").Value = 4100.4546
$ws.Range("J89").Value = 6000
$ws.Range("K89").Value = 20502.273
$ws.Range("L89").Value = 30000
$ws.Range("M89").Value = -14886.273
$ws.Range("N89").Value = -41232
$ws.Range("H91").Value = 6000
$ws.Range("J91").Value = 5000
$ws.Range("L91").Value = 5000
$ws.Range("N91").Value = -7808
$ws.Range("H99").Value = 4462
$ws.Range("I99").Value = 4449.143
$ws.Range("J99").Value = 4507
$ws.Range("K99").Value = 4449.143
$ws.Range("L99").Value = 4507
$ws.Range("M99").Value = -2951.143
$ws.Range("N99").Value = -7503
$ws.Range("H126").Value = 4462
$ws.Range("I126").Value = 4449.143
$ws.Range("J126").Value = 4507
$ws.Range("K126").Value = 13347.429
$ws.Range("L126").Value = 13521
$ws.Range("M126").Value = -10877.429
$ws.Range("N126").Value = -18461
$ws.Range("H128").Value = 30000
$ws.Range("J128").Value = 30000
$ws.Range("L128").Value = 30000
$ws.Range("N128").Value = -39960

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H113").Value = 1298.8
$ws.Range("I113").Value = 751.5
$ws.Range("K113").Value = 2254.5
$ws.Range("M113").Value = -84.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 33038
$ws.Range("I26").Value = 33038
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 33038
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -32758
$ws.Range("N26").ClearContents()
$ws.Range("H50").Value = 33038
$ws.Range("I50").Value = 33038
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 33038
$ws.Range("L50").Value = 0
$ws.Range("M50").Value = -32540
$ws.Range("N50").ClearContents()
$ws.Range("H58").Value = 37166.332
$ws.Range("J58").Value = 36999.5
$ws.Range("L58").Value = 36999.5
$ws.Range("N58").Value = -37553.5
$ws.Range("H128").Value = 39990
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 750.7059
$ws.Range("I55").Value = 619.36365
$ws.Range("K55").Value = 619.36365
$ws.Range("M55").Value = -446.36365
$ws.Range("H100").Value = 2766.3333
$ws.Range("I100").Value = 2766.3333
$ws.Range("K100").Value = 2766.3333
$ws.Range("M100").Value = -2225.3333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 54500
$ws.Range("J70").Value = 54500
$ws.Range("L70").Value = 54500
$ws.Range("N70").Value = -55130
$ws.Range("H73").Value = 54500
$ws.Range("J73").Value = 54500
$ws.Range("L73").Value = 54500
$ws.Range("N73").Value = -56684
$ws.Range("H98").Value = 29006.5
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 29006.5
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 29006.5
$ws.Range("M98").ClearContents()
$ws.Range("N98").Value = -34996.5
$ws.Range("H113").Value = 983.1429000000001
$ws.Range("I113").Value = 997
$ws.Range("J113").Value = 900
$ws.Range("K113").Value = 2991
$ws.Range("L113").Value = 2700
$ws.Range("M113").Value = -821
$ws.Range("N113").Value = -7040
$ws.Range("H122").Value = 2770.5715
$ws.Range("I122").Value = 2872.7896
$ws.Range("K122").Value = 8618.3688
$ws.Range("M122").Value = -6168.3688
